$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "95.185.27"
$ws.Range("E2").Value = "  -0.36%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.547.98"
$ws.Range("E3").Value = "  +0.11%  "

$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.71"
$ws.Range("E5").Value = "  -1.00%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "649.36"
$ws.Range("E6").Value = "  +2.36%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.47"
$ws.Range("E7").Value = "  -0.08%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.397"
$ws.Range("E8").Value = "  -0.23%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  -0.08%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.997"
$ws.Range("E10").Value = "  -1.66%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.549.15"
$ws.Range("E11").Value = "  +0.24%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.202"
$ws.Range("E12").Value = "  +1.00%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "42.18"
$ws.Range("E13").Value = "  -2.00%  "

$ws.Range("E14").Value = "  +0.86%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.229.99"
$ws.Range("E15").Value = "  +0.14%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "95.081.46"
$ws.Range("E16").Value = "  -0.40%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000252"
$ws.Range("E17").Value = "  -0.01%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.549.11"
$ws.Range("E18").Value = "  +0.12%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.89"
$ws.Range("E19").Value = "  -0.91%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.59"
$ws.Range("E20").Value = "  -2.88%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.71"
$ws.Range("E21").Value = "  -0.27%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.48"
$ws.Range("E22").Value = "  +2.48%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "503.06"
$ws.Range("E23").Value = "  -2.18%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.473"
$ws.Range("E24").Value = "  -5.19%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000194"
$ws.Range("E25").Value = "  +0.94%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.66"
$ws.Range("E26").Value = "  -0.34%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "91.61"
$ws.Range("E27").Value = "  -1.12%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.735.57"
$ws.Range("E28").Value = "  -0.12%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "12.35"
$ws.Range("E29").Value = "  +0.96%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.00"
$ws.Range("E30").Value = "  -1.39%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "11.35"
$ws.Range("E31").Value = "  -0.72%  "

$ws.Range("E32").Value = "  +0.03%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.140"
$ws.Range("E33").Value = "  -3.78%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("E34").Value = "  -0.93%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.177"
$ws.Range("E35").Value = "  -3.23%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "31.61"
$ws.Range("E36").Value = "  +5.54%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.555"
$ws.Range("E37").Value = "  -1.26%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.38"
$ws.Range("E38").Value = "  +8.48%  "

$ws.Range("E39").Value = "  +8.47%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "569.72"
$ws.Range("E40").Value = "  -1.59%  "

$ws.Range("E41").Value = "  +0.03%  "

$ws.Range("E42").Value = "  -0.64%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.896"
$ws.Range("E43").Value = "  -2.46%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.74"
$ws.Range("E44").Value = "  +0.16%  "

$ws.Range("E45").Value = "  +5.55%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "23.40"
$ws.Range("E46").Value = "  -1.83%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.60"
$ws.Range("E47").Value = "  +0.37%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "33.46"
$ws.Range("E48").Value = "  +32.22%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0409"
$ws.Range("E49").Value = "  -4.04%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.57"
$ws.Range("E50").Value = "  +0.67%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "53.09"
$ws.Range("E51").Value = "  -1.24%  "
